$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 / Row 5: E column becomes a literal phone number, F column
# becomes a hyperlink cell showing "Tcs@1983" (replacing the old "s3"/"p3"
# shared strings which drop out of sharedStrings.xml as a result). ---
$ws.Range("E4").Value = 9967887510
$ws.Range("E5").Value = 9967887510

$ws.Hyperlinks.Add($ws.Range("F4"), "https://example.com", "", "", "Tcs@1983")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://example.com", "", "", "Tcs@1983")

# Hyperlinks.Add stamps a brand-new cell style; put F4/F5 back onto the
# same "Hyperlink-ish" style already used by D4/D5 (copy/paste formats
# reuses the existing cellXfs entry instead of minting another one).
$ws.Range("D4").Copy()
$ws.Range("F4").PasteSpecial(-4122)
$ws.Range("D4").Copy()
$ws.Range("F5").PasteSpecial(-4122)
$ws.Range("A1").Copy()

# --- Selection moved to E5:F5 ---
$ws.Range("E5:F5").Select()

# --- Column E width 10.42578125 -> 11 (ColumnWidth setter has a fixed
# +5/6 char padding quirk in this host, so back it out to land on 11). ---
$ws.Columns("E").ColumnWidth = 11 - 5/6
